# edit.ps1
# Applies the "fixed results for mitigation" change: converts raw counts in
# columns C, D, F, and the language columns (G..AB) for rows 2-17 into the
# corrected percentage values, as captured by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 46.03
    "D2" = 22.64
    "F2" = 14.21
    "H2" = 0.01
    "C3" = 76.5
    "D3" = 74.67
    "F3" = 27.94
    "I3" = 0.39
    "K3" = 0.13
    "M3" = 0.13
    "O3" = 0.39
    "S3" = 0.13
    "C4" = 52.42
    "D4" = 53.3
    "F4" = 51.66
    "J4" = 0.16
    "S4" = 0.03
    "C5" = 66.73
    "D5" = 0.89
    "F5" = 18.26
    "N5" = 0.03
    "C6" = 47.7
    "D6" = 12.04
    "F6" = 17.57
    "J6" = 0.01
    "Q6" = 0.01
    "S6" = 0.02
    "U6" = 0.03
    "Y6" = 0.02
    "C7" = 77.11
    "D7" = 22.74
    "F7" = 50.44
    "I7" = 0.15
    "J7" = 0.15
    "O7" = 0.44
    "P7" = 0.15
    "Z7" = 0.29
    "AA7" = 0.29
    "C8" = 79.34999999999999
    "D8" = 13.53
    "I8" = 0.08
    "Z8" = 0.08
    "C9" = 47.44
    "D9" = 15.58
    "F9" = 27.09
    "H9" = 0.02
    "S9" = 0
    "C10" = 57.14
    "D10" = 16.43
    "F10" = 1.54
    "H10" = 0.03
    "I10" = 0.02
    "J10" = 0.01
    "S10" = 0.01
    "Y10" = 0.01
    "C11" = 53.76
    "D11" = 16.22
    "F11" = 37.78
    "G11" = 0.01
    "J11" = 0.01
    "N11" = 0.02
    "Y11" = 0.01
    "C12" = 83.72
    "D12" = 56.98
    "F12" = 9.75
    "J12" = 0.06
    "L12" = 0.1
    "S12" = 0.08
    "AA12" = 0.02
    "C13" = 53.5
    "D13" = 15.32
    "F13" = 12.9
    "J13" = 0.03
    "S13" = 0.01
    "Y13" = 0.04
    "C14" = 41.72
    "D14" = 12.56
    "F14" = 9.460000000000001
    "J14" = 0.01
    "Q14" = 0.01
    "C15" = 72.39
    "D15" = 16.08
    "F15" = 1.12
    "J15" = 0.03
    "L15" = 0.01
    "N15" = 0.01
    "Y15" = 0.03
    "C16" = 52.65
    "D16" = 12.9
    "F16" = 45.36
    "G16" = 0.04
    "I16" = 0.15
    "L16" = 0.04
    "S16" = 0.04
    "Y16" = 0.04
    "AA16" = 0.11
    "C17" = 86.31999999999999
    "D17" = 23.55
    "F17" = 3.17
    "H17" = 0.09
    "P17" = 0.27
    "T17" = 0.09
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
